$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.47434401512146
$ws.Range("B1").Value = 1.450868487358093
$ws.Range("C1").Value = 6.735958576202393
$ws.Range("D1").Value = 1.648200631141663
$ws.Range("E1").Value = 0.9800604581832886
